$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 448; this shifts the existing rows
# 448-532 down to 451-535 (dimension grows from A1:T532 to A1:T535).
$ws.Rows("448:450").Insert()

# --- New row 448: 1a amarillo, week of 44637 ---
$ws.Range("A448").Value = 2
$ws.Range("B448").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C448").Value = "Coquimbo"
$ws.Range("D448").Value2 = 44637
$ws.Range("E448").Value = 4
$ws.Range("F448").Value = "Fruta"
$ws.Range("G448").Value = 100102
$ws.Range("H448").Value = "Cítricos"
$ws.Range("I448").Value = 100102003
$ws.Range("J448").Value = "Limón"
$ws.Range("K448").Value = "Sin especificar"
$ws.Range("L448").Value = "1a amarillo"
$ws.Range("M448").Value = 750
$ws.Range("N448").Value = 17800
$ws.Range("O448").Value = 18000
$ws.Range("P448").Value = 17900
$ws.Range("Q448").Value = "$/malla 16 kilos"
$ws.Range("R448").Value = "Provincia de Limarí"
$ws.Range("S448").Value = 1119
$ws.Range("T448").Value = 16

# --- New row 449: 2a amarillo, week of 44637 ---
$ws.Range("A449").Value = 2
$ws.Range("B449").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C449").Value = "Coquimbo"
$ws.Range("D449").Value2 = 44637
$ws.Range("E449").Value = 4
$ws.Range("F449").Value = "Fruta"
$ws.Range("G449").Value = 100102
$ws.Range("H449").Value = "Cítricos"
$ws.Range("I449").Value = 100102003
$ws.Range("J449").Value = "Limón"
$ws.Range("K449").Value = "Sin especificar"
$ws.Range("L449").Value = "2a amarillo"
$ws.Range("M449").Value = 600
$ws.Range("N449").Value = 14800
$ws.Range("O449").Value = 15000
$ws.Range("P449").Value = 14900
$ws.Range("Q449").Value = "$/malla 16 kilos"
$ws.Range("R449").Value = "Provincia de Limarí"
$ws.Range("S449").Value = 931
$ws.Range("T449").Value = 16

# --- New row 450: 3a amarillo, week of 44637 ---
$ws.Range("A450").Value = 2
$ws.Range("B450").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C450").Value = "Coquimbo"
$ws.Range("D450").Value2 = 44637
$ws.Range("E450").Value = 4
$ws.Range("F450").Value = "Fruta"
$ws.Range("G450").Value = 100102
$ws.Range("H450").Value = "Cítricos"
$ws.Range("I450").Value = 100102003
$ws.Range("J450").Value = "Limón"
$ws.Range("K450").Value = "Sin especificar"
$ws.Range("L450").Value = "3a amarillo"
$ws.Range("M450").Value = 540
$ws.Range("N450").Value = 12800
$ws.Range("O450").Value = 13000
$ws.Range("P450").Value = 12900
$ws.Range("Q450").Value = "$/malla 16 kilos"
$ws.Range("R450").Value = "Provincia de Limarí"
$ws.Range("S450").Value = 806
$ws.Range("T450").Value = 16
